$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at row 100, shifting existing rows (100..231) down to (101..232).
$ws.Rows.Item(100).Insert()

# Populate the newly inserted row with the "other" -> "other" mapping.
$ws.Cells.Item(100, 1).Value = "other"
$ws.Cells.Item(100, 2).Value = "other"

# Update the view state to match (active selection, then scroll the
# viewport so the new row's neighborhood is visible).
$ws.Activate()
$ws.Range("A100").Select()
$excel.ActiveWindow.ScrollRow = 84
